$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Config")
$wsTests  = $wb.Worksheets.Item("Test Cases")

# --- Config sheet updates ---
# RunInParallel: Yes -> No
$wsConfig.Range("B3").Value = "No"
# NumberOfBrowsers: 3 -> 1
$wsConfig.Range("B4").Value = 1

# --- Test Cases sheet updates ---
# Insert a new "Description" column after TestCaseNumber (column A)
$wsTests.Range("B1").EntireColumn.Insert()

$wsTests.Range("B1").Value = "Description"
$wsTests.Range("B2").Value = "Description1"
$wsTests.Range("B3").Value = "Description2"
$wsTests.Range("B4").Value = "Description3"
$wsTests.Range("B5").Value = "Description4"
$wsTests.Range("B6").Value = "Description5"
$wsTests.Range("B7").Value = "Description6"
$wsTests.Range("B8").Value = "Description7"
$wsTests.Range("B9").Value = "Description8"
$wsTests.Range("B10").Value = "Description9"
$wsTests.Range("B11").Value = "Description10"
$wsTests.Range("B12").Value = "Description11"

# Execute (column E, after the new Description column) for row 2: Groups=Smoke -> TestCaseNumber=101,102,103
$wsTests.Range("E2").Value = "TestCaseNumber=101,102,103"

# --- Selection / active cell updates (match final authored state) ---
$wsConfig.Range("B5").Select()
$wsTests.Range("C12").Select()
